$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("peds")

function Set-CellText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value2 = $text
    $rng.Font.Name = "Arial"
    $rng.Font.Size = 10
}

Set-CellText "W15" "Catch Up Age 4 months to 18 years"
Set-CellText "X15" "Dose 1 to Dose 2: 8 weeks"
Set-CellText "W13" "Catch Up Age 4 months to 18 years"
Set-CellText "X13" "Dose 1 to Dose 2: 6 months"
Set-CellText "W12" "Catch Up Age 4 months to 6 years"
Set-CellText "X12" "3 months"
Set-CellText "Y12" "Catch Up Age 7 years to 18 years"
Set-CellText "Z12" "3 months if younger than 13 years or 4 weeks if 13 years or older"
Set-CellText "W8" "Catch Up Age 4 months to 4 years"
Set-CellText "X8" "Dose 1 to Dose 2: 4 weeks, Dose 2 to Dose 3: 4 weeks, Dose 3 to Dose 4: 6 months (minimum age is 4 years old)"
Set-CellText "Y8" "Catch Up Age 4 years to 6 years"
Set-CellText "Z8" "Dose 1 to Dose 2: 4 weeks, Dose 2 to Dose 3: 6 months"
Set-CellText "AA8" "Catch Up Age 7 years to 18 years"
Set-CellText "AB8" "Dose 1 to Dose 2: 4 weeks, Dose 2 to Dose 3: 6 months, Dose 3 to Dose 4 (If all 3 doses were administered at < 4 years or if the third dose was administered < 6 months after the second dose): 6 months"
Set-CellText "W2" "Catch Up Age 4 months to 18 years"
Set-CellText "X2" "Dose 1 to Dose 2: 4 weeks, Dose 2 to Dose 3: 8 weeks (at least 16 weeks after dose 1)"
Set-CellText "X4" "Dose 1 to Dose 2: 4 weeks, Dose 2 to Dose 3: 4 weeks, Dose 3 to Dose 4: 6 months, Dose 4 to Dose 5: 6 months"
Set-CellText "W4" "Catch Up Age  4 months to 6 years"
Set-CellText "W16" "Catch Up Age 4 months to 6 years"
Set-CellText "X16" "Dose 1 to Dose 2: 4 weeks, Dose 2 to Dose 3: 4 weeks (max age for dose 2 is 8 months)"
Set-CellText "W5" "Catch Up Age 4 months to 6 years: If 1st dose was after 15 months"
Set-CellText "X5" "No further doses needed"
Set-CellText "Y5" "Catch Up Age 4 months to 6 years: If 1st dose was before 12 months"
Set-CellText "Z5" "Dose 1 to Dose 2: 4 weeks. If current age is less than 12 months and first dose was before 7 months and 1st dose was ActHib, Pentacel, Hiberix, Vaxelis, or unknown, Dose 2 to Dose 3: 4 weeks"
Set-CellText "AA5" "Catch Up Age 4 months to 6 years: If 1st dose was between 12 and 14 months"
Set-CellText "AB5" "Dose 1 to Dose 2: 8 weeks. If current age is less than 12 months and first dose was in-between 7 to 11 months OR current age is 12-59 months and first dose was administered before the 1st birthday and second dose was administered at younger than 15 months OR if both doses of PedvaxHIB were administered before the 1st birthday, Dose 2 to Dose 3: 8 weeks (must be 12-59 months old). If all 3 doses were administered before the 1st birthday and the child is 12-59 months old, Dose 3 to Dose 4: 8 weeks"

# Update sheet view: topLeftCell and selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("T15").Select() | Out-Null
